# Quarterly indexing esoteric bug-fix operation
#
# The per-quarter error-metric table (rows 2-11, one row per quarter
# Q0..Q9) had an off-by-one indexing bug: each quarter's metrics were
# being written one row below where they belonged, so the freshly
# computed Q0 figures were dropped and every other quarter showed the
# previous quarter's numbers. This fix shifts the existing metrics
# (columns B:G) down by one row and fills row 2 (Q0) with the newly
# computed values for that quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 10,6
$arr[0,0] = -0.02907897629796788
$arr[0,1] = 0.3131278957257717
$arr[0,2] = 0.181524606355785
$arr[0,3] = 0.4260570458938391
$arr[0,4] = 0.43998257208981
$arr[0,5] = 15
$arr[1,0] = 0.0370360702193368
$arr[1,1] = 0.2877205798324301
$arr[1,2] = 0.135345957066826
$arr[1,3] = 0.3678939481247633
$arr[1,4] = 0.3798420736582994
$arr[1,5] = 14
$arr[2,0] = 0.01989109342689383
$arr[2,1] = 0.3298745939743749
$arr[2,2] = 0.1503185731831121
$arr[2,3] = 0.3877093926939507
$arr[2,4] = 0.4030092977225764
$arr[2,5] = 13
$arr[3,0] = 0.06273519545459039
$arr[3,1] = 0.2460334154802962
$arr[3,2] = 0.09558537405740077
$arr[3,3] = 0.3091688439306276
$arr[3,4] = 0.3161984462764572
$arr[3,5] = 12
$arr[4,0] = 0.02711121047383215
$arr[4,1] = 0.2779584568479502
$arr[4,2] = 0.1216586105046287
$arr[4,3] = 0.34879594393374
$arr[4,4] = 0.3647135205176791
$arr[4,5] = 11
$arr[5,0] = -0.01088135635153479
$arr[5,1] = 0.2849421119723689
$arr[5,2] = 0.09825224066429059
$arr[5,3] = 0.3134521345664926
$arr[5,4] = 0.3302084135617004
$arr[5,5] = 10
$arr[6,0] = -0.04428645741563344
$arr[6,1] = 0.3765628100937468
$arr[6,2] = 0.1690784592812682
$arr[6,3] = 0.4111915116843588
$arr[6,4] = 0.4335975266214011
$arr[6,5] = 9
$arr[7,0] = -0.01393931246739222
$arr[7,1] = 0.3509192590318558
$arr[7,2] = 0.1607817096315818
$arr[7,3] = 0.4009759464501354
$arr[7,4] = 0.4284022211487281
$arr[7,5] = 8
$arr[8,0] = -0.03688841855209302
$arr[8,1] = 0.2799590153621541
$arr[8,2] = 0.1329006694222477
$arr[8,3] = 0.3645554408073589
$arr[8,4] = 0.3917438359423185
$arr[8,5] = 7
$arr[9,0] = -0.06364682135181432
$arr[9,1] = 0.1670219060428917
$arr[9,2] = 0.03807752624502202
$arr[9,3] = 0.1951346362002964
$arr[9,4] = 0.2020691219662191
$arr[9,5] = 6

$ws.Range("B2:G11").Value = $arr
